$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gross Yield (CF)")

$ws.Range("B2").Value = 350420000
$ws.Range("B3").Value = 997442988.46
$ws.Range("B8").Value = 0.14690309
$ws.Range("B17").Value = 17498999.8
$ws.Range("D25").Value = 0

$wb.Save()
